$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / non-numeric-looking updates
$ws.Range('D2').Value = '29.351.00'
$ws.Range('E2').Value = '  +1.44%  '
$ws.Range('D3').Value = '1.948.25'
$ws.Range('E3').Value = '  +2.93%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('E5').Value = '  +0.44%  '
$ws.Range('E6').Value = '  -0.34%  '
$ws.Range('E7').Value = '  +1.11%  '
$ws.Range('E8').Value = '  +0.37%  '
$ws.Range('E9').Value = '  +0.72%  '
$ws.Range('E10').Value = '  +1.11%  '
$ws.Range('E11').Value = '  +2.08%  '
$ws.Range('D12').Value = '1.946.67'
$ws.Range('E12').Value = '  +6.83%  '
$ws.Range('E13').Value = '  +2.74%  '
$ws.Range('E14').Value = '  +1.64%  '
$ws.Range('E15').Value = '  +1.89%  '
$ws.Range('E16').Value = '  +0.13%  '
$ws.Range('E17').Value = '  -0.30%  '
$ws.Range('E18').Value = '  +0.16%  '
$ws.Range('E19').Value = '  +1.29%  '
$ws.Range('E20').Value = '  -0.27%  '
$ws.Range('D21').Value = '29.397.03'
$ws.Range('E21').Value = '  +1.68%  '
$ws.Range('E22').Value = '  +4.59%  '
$ws.Range('E23').Value = '  +2.02%  '
$ws.Range('D24').Value = '2.160.49'
$ws.Range('E24').Value = '  +4.26%  '
$ws.Range('E25').Value = '  +2.12%  '
$ws.Range('E26').Value = '  +0.33%  '
$ws.Range('E27').Value = '  +1.39%  '
$ws.Range('E28').Value = '  -0.16%  '
$ws.Range('E29').Value = '  +1.25%  '
$ws.Range('E30').Value = '  -1.92%  '
$ws.Range('E31').Value = '  +0.11%  '
$ws.Range('E32').Value = '  -1.40%  '
$ws.Range('E33').Value = '  -1.17%  '
$ws.Range('E34').Value = '  -0.12%  '
$ws.Range('B35').Value = 'PEPE'
$ws.Range('C35').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('E35').Value = '  +139.19%  '
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('E36').Value = '  -3.51%  '
$ws.Range('E37').Value = '  +0.71%  '
$ws.Range('E38').Value = '  -1.26%  '
$ws.Range('E39').Value = '  +2.22%  '
$ws.Range('E40').Value = '  -0.47%  '
$ws.Range('E41').Value = '  +0.00%  '
$ws.Range('E42').Value = '  +0.88%  '
$ws.Range('E43').Value = '  +2.83%  '
$ws.Range('E44').Value = '  +0.38%  '
$ws.Range('E45').Value = '  +0.66%  '
$ws.Range('E46').Value = '  -1.73%  '
$ws.Range('E47').Value = '  +0.06%  '
$ws.Range('E48').Value = '  -1.29%  '
$ws.Range('E49').Value = '  +2.82%  '
$ws.Range('E50').Value = '  +0.97%  '
$ws.Range('E51').Value = '  +1.15%  '

# Numeric-looking values that must remain text: force Text format, set, then restore default style
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '326.82'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9980'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4629'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3919'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07895'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.000'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '22.38'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.861'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.133'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.07066'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '88.09'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.9996'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000009978'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.17'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9988'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.530'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.21'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.101'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '156.43'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.53'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.939'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '119.12'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.894'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09355'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.8997'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.237'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.332'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.000003856'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.163'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.05806'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.176'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.02116'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9963'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.726'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5737'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1821'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '9.808'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '11.97'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.222'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5357'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.06952'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.601'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.860'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '113.85'
$ws.Range('D51').Style = 'Normal'
